# Adding to the diary #1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team name (B1), plain/default style ---
$ws.Range("B1").Value = "acdc"

# --- Group members (B2, C2) -- small Arial 10pt font ---
$ws.Range("B2").Value = "mpha0054"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").Font.Name = "Arial"

$ws.Range("C2").Value = "kngu0086"
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# --- New meeting entry, row 10 (copy formatting from the existing row 7 entries) ---
$ws.Range("A10").Value = 45186
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B10").Value = 0.48958333333333331
$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C10").Value = 0.59722222222222221
$ws.Range("C7").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("D10").Value = "All"

$ws.Range("E10").Value = "Choosing the topic for the assignment and deciding on workflow."
$ws.Range("E7").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Row heights (wrap-text rows re-flow on recalculation)
$ws.Rows.Item(7).RowHeight = 78
$ws.Rows.Item(10).RowHeight = 31.2

# Printable page orientation
$ws.PageSetup.Orientation = 1

# Restore active selection
$ws.Range("D12").Select()
